$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before) values for the columns that rotate between rows 2, 3 and 4
$D2 = $ws.Range("D2").Value2
$D3 = $ws.Range("D3").Value2
$D4 = $ws.Range("D4").Value2

$M2 = $ws.Range("M2").Value2
$M3 = $ws.Range("M3").Value2
$M4 = $ws.Range("M4").Value2

$O2 = $ws.Range("O2").Value2
$O3 = $ws.Range("O3").Value2
$O4 = $ws.Range("O4").Value2

$P2 = $ws.Range("P2").Value2
$P3 = $ws.Range("P3").Value2
$P4 = $ws.Range("P4").Value2

$S2 = $ws.Range("S2").Value2
$S3 = $ws.Range("S3").Value2
$S4 = $ws.Range("S4").Value2

# Rotate: row2 <- old row4, row3 <- old row2, row4 <- old row3
$ws.Range("D2").Value2 = $D4
$ws.Range("D3").Value2 = $D2
$ws.Range("D4").Value2 = $D3

$ws.Range("M2").Value2 = $M4
$ws.Range("M3").Value2 = $M2
$ws.Range("M4").Value2 = $M3

$ws.Range("O2").Value2 = $O4
$ws.Range("O3").Value2 = $O2
$ws.Range("O4").Value2 = $O3

$ws.Range("P2").Value2 = $P4
$ws.Range("P3").Value2 = $P2
$ws.Range("P4").Value2 = $P3

$ws.Range("S2").Value2 = $S4
$ws.Range("S3").Value2 = $S2
$ws.Range("S4").Value2 = $S3
